{"js": "// Fill in the previously-empty \"Simulation\" test-case row of the third\n// table (the \"Testen\" / testing-effort table) with the new entries that\n// the author added: task = \"Simulation\", name = \"Kalauner Paul\",\n// time = \"35 Minuten\".\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// The \"Testen\" table is the 3rd table in the document (index 2): it lists\n// \"Klassen Lager und Lagermitarbeiter\" / \"Klassen Lieferant, ...\" followed\n// by two still-empty rows. We fill the first of those empty rows.\nconst testTable = tables.items[2];\ntestTable.rows.load(\"items\");\nawait context.sync();\n\nconst targetRow = testTable.rows.items[2];\ntargetRow.cells.load(\"items\");\nawait context.sync();\n\nconst values = [\"Simulation\", \"Kalauner Paul\", \"35 Minuten\"];\nfor (let i = 0; i < targetRow.cells.items.length && i < values.length; i++) {\n  targetRow.cells.items[i].body.insertText(values[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Fill in the previously-empty \"Simulation\" test-case row of the third\n# table (the \"Testen\" / testing-effort table) with the new entries that\n# the author added: task = \"Simulation\", name = \"Kalauner Paul\",\n# time = \"35 Minuten\".\n$d = $word.ActiveDocument\n\n# The \"Testen\" table is the 3rd table in the document: it lists\n# \"Klassen Lager und Lagermitarbeiter\" / \"Klassen Lieferant, ...\" followed\n# by two still-empty rows. We fill the first of those empty rows (row 3).\n$table = $d.Tables.Item(3)\n\n$table.Cell(3, 1).Range.Text = \"Simulation\"\n$table.Cell(3, 2).Range.Text = \"Kalauner Paul\"\n$table.Cell(3, 3).Range.Text = \"35 Minuten\"\n"}
